# Updated cryptos list on Fri Mar  8 19:17:03 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are plain-text "numbers" in the source data
# (thousand separators use '.', same as the decimal separator) - force
# text storage so Excel's COM layer doesn't silently coerce them into
# real numeric values.
$priceCells = "D2","D3","D4","D5","D6","D9","D11","D12","D14","D15","D16","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D35","D36","D38","D40","D43","D44","D47","D48","D50","D51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.618.90"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.919.70"
$ws.Range("E3").Value = "  +1.13%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "485.02"
$ws.Range("E5").Value = "  +4.36%  "

# Row 6 - Solana
$ws.Range("D6").Value = "147.40"
$ws.Range("E6").Value = "  -0.43%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.99%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").Value = "  -3.91%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +7.22%  "

# Row 11 - ShibaInu
$ws.Range("D11").Value = "0.0000355"
$ws.Range("E11").Value = "  +12.26%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "42.45"
$ws.Range("E12").Value = "  -3.77%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +0.49%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.551.65"
$ws.Range("E14").Value = "  +1.12%  "

# Row 15 - Uniswap
$ws.Range("D15").Value = "14.58"
$ws.Range("E15").Value = "  -1.33%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.923.55"
$ws.Range("E16").Value = "  +1.04%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.39%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "19.74"
$ws.Range("E18").Value = "  -1.98%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  -3.03%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "68.608.98"
$ws.Range("E20").Value = "  +1.24%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "432.07"
$ws.Range("E21").Value = "  -0.21%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "14.53"
$ws.Range("E22").Value = "  -2.44%  "

# Row 23 - ImmutableX
$ws.Range("D23").Value = "3.33"
$ws.Range("E23").Value = "  +1.44%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "86.97"
$ws.Range("E24").Value = "  -1.63%  "

# Row 25 - Filecoin
$ws.Range("D25").Value = "11.34"
$ws.Range("E25").Value = "  +10.33%  "

# Row 26 - was RenderToken, now PancakeSwap
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "3.59"
$ws.Range("E26").Value = "  +0.55%  "

# Row 27 - was PancakeSwap, now RenderToken
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.66"
$ws.Range("E27").Value = "  +3.14%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "38.09"
$ws.Range("E28").Value = "  +0.80%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  +6.74%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "715.49"
$ws.Range("E30").Value = "  -4.28%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "13.22"
$ws.Range("E31").Value = "  -4.78%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -4.18%  "

# Row 34 - PEPE
$ws.Range("E34").Value = "  +32.35%  "

# Row 35 - InjectiveProtocol
$ws.Range("D35").Value = "41.37"
$ws.Range("E35").Value = "  -4.73%  "

# Row 36 - OKB
$ws.Range("D36").Value = "58.45"
$ws.Range("E36").Value = "  +1.82%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  -7.70%  "

# Row 38 - NEARProtocol
$ws.Range("D38").Value = "5.56"
$ws.Range("E38").Value = "  -0.99%  "

# Row 39 - Dai
$ws.Range("E39").Value = "  -0.15%  "

# Row 40 - Fetch.AI
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +8.71%  "

# Row 41 - VeChain
$ws.Range("E41").Value = "  -2.42%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +10.56%  "

# Row 43 - ThetaToken
$ws.Range("D43").Value = "2.98"
$ws.Range("E43").Value = "  +1.58%  "

# Row 44 - TheGraph
$ws.Range("D44").Value = "0.344"
$ws.Range("E44").Value = "  -1.90%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  -1.39%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  -0.02%  "

# Row 47 - was ARBITRUM, now LidoDAOToken
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").Value = "3.40"
$ws.Range("E47").Value = "  -1.28%  "

# Row 48 - was LidoDAOToken, now ARBITRUM
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "2.16"
$ws.Range("E48").Value = "  +1.14%  "

# Row 49 - ApeXProtocol
$ws.Range("E49").Value = "  -2.00%  "

# Row 50 - Monero
$ws.Range("D50").Value = "147.82"
$ws.Range("E50").Value = "  +1.76%  "

# Row 51 - Stacks
$ws.Range("D51").Value = "2.83"
$ws.Range("E51").Value = "  -2.62%  "
